$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on price cells that would otherwise be
# auto-converted to numeric values by Excel (single-decimal-point values).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '26.310.08'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '1.596.43'
$ws.Range("E3").Value = '  +0.48%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '211.59'
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("D10").Value = '19.04'
$ws.Range("E11").Value = '  +1.23%  '
$ws.Range("E12").Value = '  +0.46%  '
$ws.Range("D13").Value = '1.582.66'
$ws.Range("E13").Value = '  -0.38%  '
$ws.Range("E14").Value = '  -0.44%  '
$ws.Range("D15").Value = '0.504'
$ws.Range("E15").Value = '  -1.10%  '
$ws.Range("D16").Value = '63.47'
$ws.Range("E16").Value = '  -0.11%  '
$ws.Range("D17").Value = '26.302.81'
$ws.Range("E17").Value = '  +0.48%  '
$ws.Range("D18").Value = '230.62'
$ws.Range("E18").Value = '  +7.77%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.0₃0721'
$ws.Range("E19").Value = '  -0.18%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '7.62'
$ws.Range("E20").Value = '  +3.86%  '
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("E23").Value = '  +2.45%  '
$ws.Range("D25").Value = '146.67'
$ws.Range("E25").Value = '  +1.40%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("E27").Value = '  +0.33%  '
$ws.Range("E28").Value = '  +0.36%  '
$ws.Range("D29").Value = '15.37'
$ws.Range("E29").Value = '  +2.11%  '
$ws.Range("E31").Value = '  +0.21%  '
$ws.Range("D32").Value = '1.514.87'
$ws.Range("E32").Value = '  +6.76%  '
$ws.Range("E33").Value = '  +1.44%  '
$ws.Range("E34").Value = '  -0.48%  '
$ws.Range("E35").Value = '  -0.32%  '
$ws.Range("E36").Value = '  +0.63%  '
$ws.Range("E37").Value = '  -2.93%  '
$ws.Range("E38").Value = '  -0.37%  '
$ws.Range("D39").Value = '0.816'
$ws.Range("E39").Value = '  -0.69%  '
$ws.Range("E40").Value = '  -1.46%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").Value = '2.17'
$ws.Range("E42").Value = '  +1.84%  '
$ws.Range("E43").Value = '  -3.12%  '
$ws.Range("D44").Value = '1.734.02'
$ws.Range("E44").Value = '  +0.51%  '
$ws.Range("E45").Value = '  -0.58%  '
$ws.Range("D46").Value = '60.48'
$ws.Range("E46").Value = '  -0.77%  '
$ws.Range("D47").Value = '88.36'
$ws.Range("E47").Value = '  +1.71%  '
$ws.Range("E48").Value = '  +0.29%  '
$ws.Range("E49").Value = '  -0.15%  '
$ws.Range("D50").Value = '0.0957'
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("E51").Value = '  +0.03%  '
